$d = $word.ActiveDocument

# 1. Intro paragraph: date + jurisdiction
$d.Content.Find.Execute(
    "THIS DEED OF PARTNERSHIP is executed on this 01/01/2024 at Ratnagiri, Maharashtra by and between:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "THIS DEED OF PARTNERSHIP is executed on this 2025-05-15 at [Jurisdiction Not Provided] by and between:",
    2)

# 2. Partner No. 1 name redaction
$d.Content.Find.Execute(
    "1. Advait Milind Kulkarni, Son of Milind Shashikant Kulkarni, Age 25, residing at 557/H1,Thiba Palace Road ,AnandNagar, Ratnagiri, Maharashtra (Hereinafter referred to as Partner No. 1)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. [Full name], Son of Milind Shashikant Kulkarni, Age 25, residing at 557/H1,Thiba Palace Road ,AnandNagar, Ratnagiri, Maharashtra (Hereinafter referred to as Partner No. 1)",
    2)

# 3. Partner No. 2 name redaction + new address
$d.Content.Find.Execute(
    "2. Tanmay Abhay Joshi, Son of Abhay Joshi, Age 25, residing at 557/H1,Thiba Palace Road ,AnandNagar, Ratnagiri, Maharashtra (Hereinafter referred to as Partner No. 2)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2. [Full name], Son of Abhay Joshi, Age 25, residing at Omkar Sanjiwani Apartment Ratnagiri 415612 (Hereinafter referred to as Partner No. 2)",
    2)

# 4. Place of business
$d.Content.Find.Execute(
    "The place of business shall be 557/H1,Thiba Palace Road ,AnandNagar, Ratnagiri, Maharashtra, and area of operation will be Ratnagiri.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The place of business shall be [Business address Not Provided], and area of operation will be [Area of operation Not Provided].",
    2)

# 5. Partnership commence date
$d.Content.Find.Execute(
    "The partnership shall commence on 01/01/2025 and shall be a Partnership at Will.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The partnership shall commence on 2025-05-15 and shall be a Partnership at Will.",
    2)

# 6. Duties/Governing law paragraph replacement
$d.Content.Find.Execute(
    "Advait labour and Tanmay sheth",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Governing Law: This Agreement shall be governed by and construed in accordance with the laws of the State of [insert state], without regard to its conflict of law principles. Any dispute arising under or related to this Agreement shall be subject to the exclusive jurisdiction of the courts of [insert county], [insert state].",
    2)

# 7. Signature table cell - Partner No. 1 name (this exact standalone
#    string only remains in the table after step 2 already rewrote the
#    longer intro-paragraph occurrence)
$d.Content.Find.Execute(
    "Advait Milind Kulkarni",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[Full name]",
    2)

# 8. Signature table cell - Partner No. 2 name (same reasoning as above)
$d.Content.Find.Execute(
    "Tanmay Abhay Joshi",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[Full name]",
    2)
